$wb = $excel.ActiveWorkbook

# --- "actor" sheet (sheet1): remove chart title, adjust view ---
$wsActor = $wb.Worksheets.Item("actor")
$actorChart = $wsActor.ChartObjects(1).Chart
$actorChart.HasTitle = $false

# --- "stance" sheet (sheet2): update labels, drop leftover "augmentation" row, remove chart title ---
$wsStance = $wb.Worksheets.Item("stance")
$wsStance.Range("A2").Value = "Prompt includes examples only (no definition)"
$wsStance.Range("A3").Value = "Prompt includes examples and definitions"
$wsStance.Rows.Item(4).Delete()

$stanceChart = $wsStance.ChartObjects(1).Chart
$stanceChart.HasTitle = $false

# --- view state: "stance" becomes the active/selected tab ---
$wsActor.Activate()
$excel.ActiveWindow.Zoom = 106
$wsActor.Range("C3").Select()

$wsStance.Activate()
$excel.ActiveWindow.Zoom = 109
$wsStance.Range("G5").Select()
